$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# --- Populate CO Number (col A) / Comment (col B) for rows 2-9 --------------
# Column A holds numeric-looking CO numbers that must be stored as *text*
# (they are looked up / concatenated downstream), so we force the Text
# number format before writing the value - this stops Excel's automatic
# numeric inference from turning them into number cells.
# Column B holds ordinary text ("1CO Created", "2CO Created", ...) which is
# already stored as text with no extra coercion required.

$data = @(
    @("3014010588", "1CO Created"),
    @("3014010589", "2CO Created"),
    @("3014010590", "3CO Created"),
    @("3014010591", "4CO Created"),
    @("3014010592", "5CO Created"),
    @("3014010593", "6CO Created"),
    @("3014010594", "7CO Created"),
    @("3014010595", "8CO Created")
)

$row = 2
foreach ($pair in $data) {
    $coCell = $ws.Cells.Item($row, 1)
    $coCell.NumberFormat = "@"
    $coCell.Value = $pair[0]

    $ws.Cells.Item($row, 2).Value = $pair[1]

    $row = $row + 1
}

# --- Widen column A and drop the old best-fit auto width --------------------
$ws.Columns("A").ColumnWidth = 22.14

# --- Move the active selection ----------------------------------------------
[void]$ws.Range("E15").Select()
